$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (C) column for existing data rows 2-27 from 45393 to 45394
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 3).Value2 = 45394
}

# Insert a new row above the current row 28 (the "A 11790-2024" row), shifting it to row 29
$ws.Rows.Item(28).EntireRow.Insert()

# Populate the newly inserted row 28 with the new record
$ws.Cells.Item(28, 1).Value2 = "A 10085-2024"
$ws.Cells.Item(28, 2).Value2 = 45364
$ws.Cells.Item(28, 3).Value2 = 45394
$ws.Cells.Item(28, 4).Value2 = "OKÄNT"
$ws.Cells.Item(28, 5).Value2 = "OKÄNT"
$ws.Cells.Item(28, 7).Value2 = 3.9
$ws.Cells.Item(28, 8).Value2 = 0
$ws.Cells.Item(28, 9).Value2 = 0
$ws.Cells.Item(28, 10).Value2 = 0
$ws.Cells.Item(28, 11).Value2 = 0
$ws.Cells.Item(28, 12).Value2 = 0
$ws.Cells.Item(28, 13).Value2 = 0
$ws.Cells.Item(28, 14).Value2 = 0
$ws.Cells.Item(28, 15).Value2 = 0
$ws.Cells.Item(28, 16).Value2 = 0
$ws.Cells.Item(28, 17).Value2 = 0

# Ensure the row height matches the rest of the data rows
$ws.Rows.Item(28).RowHeight = 15

# Update the "Förändrad" (C) column for the row that shifted down to 29
$ws.Cells.Item(29, 3).Value2 = 45394
